$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""
$ws.Range("H70").Value = 1466.6666
$ws.Range("I70").Value = 1550
$ws.Range("J70").Value = 1300
$ws.Range("K70").Value = 4650
$ws.Range("L70").Value = 3900
$ws.Range("M70").Value = -4380
$ws.Range("N70").Value = -4440
$ws.Range("H73").Value = 1466.6666
$ws.Range("I73").Value = 1550
$ws.Range("J73").Value = 1300
$ws.Range("K73").Value = 4650
$ws.Range("L73").Value = 3900
$ws.Range("M73").Value = -3714
$ws.Range("N73").Value = -5772
$ws.Range("H92").Value = 1184.7
$ws.Range("I92").Value = 983
$ws.Range("K92").Value = 983
$ws.Range("M92").Value = 265
$ws.Range("H103").Value = 833646.7
$ws.Range("I103").Value = 833646.7
$ws.Range("K103").Value = 2500940.1
$ws.Range("M103").Value = -2500354.1
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140
$ws.Range("H137").Value = 38406.035
$ws.Range("I137").Value = 1168.4706
$ws.Range("J137").Value = 101709.9
$ws.Range("K137").Value = 3505.4118
$ws.Range("L137").Value = 305129.7
$ws.Range("M137").Value = -955.4118000000003
$ws.Range("N137").Value = -310229.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1834.3334
$ws.Range("I102").Value = 1627.5
$ws.Range("K102").Value = 1627.5
$ws.Range("M102").Value = -5.5
$ws.Range("H114").Value = 27366.334
$ws.Range("J114").Value = 27366.334
$ws.Range("L114").Value = 27366.334
$ws.Range("N114").Value = -36044.334
$ws.Range("H135").Value = 22470.75
$ws.Range("J135").Value = 22470.75
$ws.Range("L135").Value = 22470.75
$ws.Range("N135").Value = -32610.75
$ws.Range("H139").Value = 41998
$ws.Range("J139").Value = 41998
$ws.Range("L139").Value = 41998
$ws.Range("N139").Value = -52278

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 22085.834
$ws.Range("J81").Value = 22085.834
$ws.Range("L81").Value = 22085.834
$ws.Range("N81").Value = -24207.834
$ws.Range("H84").Value = 22085.834
$ws.Range("J84").Value = 22085.834
$ws.Range("L84").Value = 66257.50199999999
$ws.Range("N84").Value = -76865.50199999999
$ws.Range("H94").Value = 854.4314000000001
$ws.Range("I94").Value = 792.6977000000001
$ws.Range("K94").Value = 792.6977000000001
$ws.Range("M94").Value = -341.6977000000001
$ws.Range("H99").Value = 1815.3334
$ws.Range("I99").Value = 1277.5
$ws.Range("K99").Value = 1277.5
$ws.Range("M99").Value = 220.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1092.2222
$ws.Range("I16").Value = 1071.6666
$ws.Range("J16").Value = 1133.3334
$ws.Range("K16").Value = 1071.6666
$ws.Range("L16").Value = 1133.3334
$ws.Range("M16").Value = -784.6666
$ws.Range("N16").Value = -1707.3334
$ws.Range("H31").Value = 14123.538
$ws.Range("I31").Value = 19184.055
$ws.Range("J31").Value = 2737.375
$ws.Range("K31").Value = 19184.055
$ws.Range("L31").Value = 2737.375
$ws.Range("M31").Value = -18889.055
$ws.Range("N31").Value = -3327.375
$ws.Range("H34").Value = 14123.538
$ws.Range("I34").Value = 19184.055
$ws.Range("J34").Value = 2737.375
$ws.Range("K34").Value = 19184.055
$ws.Range("L34").Value = 2737.375
$ws.Range("M34").Value = -18982.055
$ws.Range("N34").Value = -3141.375
$ws.Range("H113").Value = 1092.2222
$ws.Range("I113").Value = 1071.6666
$ws.Range("J113").Value = 1133.3334
$ws.Range("K113").Value = 1071.6666
$ws.Range("L113").Value = 1133.3334
$ws.Range("M113").Value = 1098.3334
$ws.Range("N113").Value = -5473.3334
$ws.Range("H135").Value = 36778.3
$ws.Range("J135").Value = 35452.668
$ws.Range("L135").Value = 35452.668
$ws.Range("N135").Value = -45592.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1107.0834
$ws.Range("I5").Value = 553.3333
$ws.Range("J5").Value = 2768.3333
$ws.Range("K5").Value = 1659.9999
$ws.Range("L5").Value = 8304.999899999999
$ws.Range("M5").Value = -1547.9999
$ws.Range("N5").Value = -8528.999899999999
$ws.Range("H63").Value = 4478.7144
$ws.Range("I63").Value = 1911.5
$ws.Range("J63").Value = 5505.6
$ws.Range("K63").Value = 5734.5
$ws.Range("L63").Value = 16516.8
$ws.Range("M63").Value = -4985.5
$ws.Range("N63").Value = -18014.8
$ws.Range("H66").Value = 4478.7144
$ws.Range("I66").Value = 1911.5
$ws.Range("J66").Value = 5505.6
$ws.Range("K66").Value = 17203.5
$ws.Range("L66").Value = 49550.4
$ws.Range("M66").Value = -13459.5
$ws.Range("N66").Value = -57038.4
$ws.Range("H131").Value = 709.41
$ws.Range("J131").Value = 710.51514
$ws.Range("L131").Value = 2131.54542
$ws.Range("N131").Value = -12211.54542
$ws.Range("H133").Value = 3837.5
$ws.Range("I133").Value = 3000
$ws.Range("J133").Value = 3957.1428
$ws.Range("K133").Value = 9000
$ws.Range("L133").Value = 11871.4284
$ws.Range("M133").Value = -3940
$ws.Range("N133").Value = -21991.4284
$ws.Range("H134").Value = 1213.56
$ws.Range("I134").Value = 1147.5
$ws.Range("J134").Value = 2799
$ws.Range("K134").Value = 3442.5
$ws.Range("L134").Value = 8397
$ws.Range("M134").Value = 1627.5
$ws.Range("N134").Value = -18537
$ws.Range("H135").Value = 1107.0834
$ws.Range("I135").Value = 553.3333
$ws.Range("J135").Value = 2768.3333
$ws.Range("K135").Value = 4979.9997
$ws.Range("L135").Value = 24914.9997
$ws.Range("M135").Value = -2444.9997
$ws.Range("N135").Value = -29984.9997
$ws.Range("H136").Value = 2643.5454
$ws.Range("I136").Value = 1761.25
$ws.Range("J136").Value = 4996.3335
$ws.Range("K136").Value = 5283.75
$ws.Range("L136").Value = 14989.0005
$ws.Range("M136").Value = -183.75
$ws.Range("N136").Value = -25189.0005
$ws.Range("H137").Value = 2939
$ws.Range("J137").Value = 4683.25
$ws.Range("L137").Value = 14049.75
$ws.Range("N137").Value = -24249.75
$ws.Range("H139").Value = 1511.2632
$ws.Range("I139").Value = 1079.625
$ws.Range("J139").Value = 3813.3333
$ws.Range("K139").Value = 3238.875
$ws.Range("L139").Value = 11439.9999
$ws.Range("M139").Value = 1901.125
$ws.Range("N139").Value = -21719.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 49490
$ws.Range("J135").Value = 48986.668
$ws.Range("L135").Value = 48986.668
$ws.Range("N135").Value = -59126.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4141.3
$ws.Range("I61").Value = 1678.6154
$ws.Range("J61").Value = 8714.857
$ws.Range("K61").Value = 1678.6154
$ws.Range("L61").Value = 8714.857
$ws.Range("M61").Value = -1476.6154
$ws.Range("N61").Value = -9118.857
$ws.Range("H68").Value = 2333.2942
$ws.Range("I68").Value = 2072.5
$ws.Range("J68").Value = 2959.2
$ws.Range("K68").Value = 2072.5
$ws.Range("L68").Value = 2959.2
$ws.Range("M68").Value = -1323.5
$ws.Range("N68").Value = -4457.2
$ws.Range("H71").Value = 2333.2942
$ws.Range("I71").Value = 2072.5
$ws.Range("J71").Value = 2959.2
$ws.Range("K71").Value = 10362.5
$ws.Range("L71").Value = 14796
$ws.Range("M71").Value = -6618.5
$ws.Range("N71").Value = -22284
$ws.Range("H113").Value = 4141.3
$ws.Range("I113").Value = 1678.6154
$ws.Range("J113").Value = 8714.857
$ws.Range("K113").Value = 1678.6154
$ws.Range("L113").Value = 8714.857
$ws.Range("M113").Value = 491.3846000000001
$ws.Range("N113").Value = -13054.857
$ws.Range("H132").Value = 1653.2188
$ws.Range("I132").Value = 1057.381
$ws.Range("K132").Value = 3172.143
$ws.Range("M132").Value = -642.143
$ws.Range("H136").Value = 28813.611
$ws.Range("I136").Value = 34296.332
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 102888.996
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -100338.996
$ws.Range("N136").Value = -9300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4198.857
$ws.Range("J62").Value = 4398.8335
$ws.Range("L62").Value = 4398.8335
$ws.Range("N62").Value = -5646.8335
$ws.Range("H65").Value = 4198.857
$ws.Range("J65").Value = 4398.8335
$ws.Range("L65").Value = 21994.1675
$ws.Range("N65").Value = -28234.1675
$ws.Range("H119").Value = 21000
$ws.Range("J119").Value = 21000
$ws.Range("L119").Value = 21000
$ws.Range("N119").Value = -30676
